$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 16.427638310189508
$ws.Range("C2").Value = 11.017066881382448
$ws.Range("D2").Value = 13.448417733397036
$ws.Range("E2").Value = -0.25453656396425117

$ws.Range("B3").Value = 34.135045502966477
$ws.Range("C3").Value = 3.4386750814914819
$ws.Range("D3").Value = 1.7820797767813872
$ws.Range("E3").Value = -0.1021508699954552

$ws.Range("B1:E3").Select()
